# C195 Remaining To Do List.docx - "fixed bugs in the system"
#
# Strategy: use Range.Find.Execute to locate precise text anchors, then
# Range.InsertXML(...) to replace exactly that range's contents with the
# desired WordprocessingML (multiple runs, proofErr markers, symbols,
# etc). InsertXML only rewrites the targeted range, leaving the rest of
# the paragraph / document (and paragraph marks outside the range)
# untouched - so it behaves like a surgical "select this text -> retype
# it with new runs/marks" edit, same as what Word's own paste/AutoCorrect
# machinery produces.
#
# Edits are applied from the bottom of the document upward so that
# earlier (still-untouched) Paragraphs(...) indices stay valid for the
# whole script.

$d = $word.ActiveDocument
$WordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-WordPkg([string]$bodyXml) {
    return @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document $WordNs><w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
}

function Get-RangeForText([string]$text) {
    $r = $d.Content
    $ok = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $text"
    }
    return $r
}

# ---------------------------------------------------------------------
# Edit A: delete "[X]  -  View Schedule..." / blank / "[X] - BUG..."
# paragraphs (20-22) entirely - their content moves into paragraph 15
# in Edit B below, so this whole span just goes away.
# ---------------------------------------------------------------------
$startP = $d.Paragraphs(20).Range
$endP = $d.Paragraphs(22).Range
$delRange = $d.Range($startP.Start, $endP.End)
$delRange.InsertXML(New-WordPkg(""))

# ---------------------------------------------------------------------
# Edit B: merge "[X]    -  Check if user has..." (15) through
# "...is coming back as null? " (18) into one paragraph reusing the
# "View Schedule by Contact ID" text (previously paragraph 20, now
# folded in here). Paragraph 19 (blank) is left alone.
# ---------------------------------------------------------------------
$p15 = $d.Paragraphs(15).Range
$p18 = $d.Paragraphs(18).Range
$mergeRange1 = $d.Range($p15.Start, $p18.End)
$body1 = @'
<w:p><w:r><w:t>[X</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>]  -</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">  View Schedule by Contact ID - Report 2  Need to figure out what to pass into initialize method</w:t></w:r></w:p>
'@
$mergeRange1.InsertXML(New-WordPkg($body1))

# ---------------------------------------------------------------------
# Edit C: title "Completed But Not Working Properly" -> split "But"
# out with gramStart/gramEnd proofing marks, same rPr on every run.
# ---------------------------------------------------------------------
$rTitle2 = Get-RangeForText("Completed But Not Working Properly")
$rPr32 = '<w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'
$body2 = @"
<w:p><w:r>$rPr32<w:t xml:space="preserve">Completed </w:t></w:r><w:proofErr w:type="gramStart"/><w:r>$rPr32<w:t>But</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>$rPr32<w:t xml:space="preserve"> Not Working Properly</w:t></w:r></w:p>
"@
$rTitle2.InsertXML(New-WordPkg($body2))

# ---------------------------------------------------------------------
# Edit D: "[  ]  -  Generate " -> split off "[  ]" with gramStart/End,
# leaving the rest of the paragraph (JavaDocs / trailing sentence)
# untouched.
# ---------------------------------------------------------------------
$rGenerate = Get-RangeForText("[  ]  -  Generate ")
$body3 = @'
<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>[  ]</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">  -  Generate </w:t></w:r></w:p>
'@
$rGenerate.InsertXML(New-WordPkg($body3))

# ---------------------------------------------------------------------
# Edit E: merge "[  ] - Discuss report..." (3) through
# "...Scene Builder" (9) into a single paragraph that now carries what
# used to be the "Create Report to view Customer Appointments..." text
# (previously paragraph 7), with the "[  ]" also split out via
# gramStart/gramEnd, and the paragraph format gains <w:u w:val="single"/>.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$p9 = $d.Paragraphs(9).Range
$mergeRange2 = $d.Range($p3.Start, $p9.End)
$body4 = @'
<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>[  ]</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">  -  Create Report to view Customer Appointments by Month </w:t></w:r><w:r><w:t>/ Type / Count</w:t></w:r><w:r><w:t xml:space="preserve"> - Report 1</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:sym w:font="Wingdings" w:char="F0DF"/></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> CURRENTLY WORKING ON</w:t></w:r></w:p>
'@
$mergeRange2.InsertXML(New-WordPkg($body4))

# ---------------------------------------------------------------------
# Edit F: title "Remaining To Do List" -> split "To" out with
# gramStart/gramEnd proofing marks, same rPr on every run.
# ---------------------------------------------------------------------
$rTitle1 = Get-RangeForText("Remaining To Do List")
$body5 = @"
<w:p><w:r>$rPr32<w:t xml:space="preserve">Remaining </w:t></w:r><w:proofErr w:type="gramStart"/><w:r>$rPr32<w:t>To</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>$rPr32<w:t xml:space="preserve"> Do List</w:t></w:r></w:p>
"@
$rTitle1.InsertXML(New-WordPkg($body5))

Write-Output "done"
